$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6170212765957447
$ws.Range("C2").Value = 0.7631578947368421
$ws.Range("D2").Value = 0.6823529411764706
$ws.Range("B3").Value = 0.8043478260869565
$ws.Range("C3").Value = 0.6727272727272727
$ws.Range("D3").Value = 0.7326732673267325
$ws.Range("B4").Value = 0.7096774193548387
$ws.Range("C4").Value = 0.7096774193548387
$ws.Range("D4").Value = 0.7096774193548387
$ws.Range("E4").Value = 0.7096774193548387
$ws.Range("B5").Value = 0.7106845513413507
$ws.Range("C5").Value = 0.7179425837320574
$ws.Range("D5").Value = 0.7075131042516016
$ws.Range("B6").Value = 0.7278057951120527
$ws.Range("C6").Value = 0.7096774193548387
$ws.Range("D6").Value = 0.7121122738459804
$ws.Range("B7").Value = 0.696969696969697
$ws.Range("C7").Value = 0.6052631578947368
$ws.Range("D7").Value = 0.6478873239436619
$ws.Range("B8").Value = 0.75
$ws.Range("C8").Value = 0.8181818181818182
$ws.Range("D8").Value = 0.7826086956521738
$ws.Range("B9").Value = 0.7311827956989247
$ws.Range("C9").Value = 0.7311827956989247
$ws.Range("D9").Value = 0.7311827956989247
$ws.Range("E9").Value = 0.7311827956989247
$ws.Range("B10").Value = 0.7234848484848485
$ws.Range("C10").Value = 0.7117224880382775
$ws.Range("D10").Value = 0.7152480097979179
$ws.Range("B11").Value = 0.7283317041381557
$ws.Range("C11").Value = 0.7311827956989247
$ws.Range("D11").Value = 0.7275612534486958
$ws.Range("B12").Value = 0.7083333333333334
$ws.Range("C12").Value = 0.4473684210526316
$ws.Range("D12").Value = 0.5483870967741936
$ws.Range("B13").Value = 0.6956521739130435
$ws.Range("C13").Value = 0.8727272727272727
$ws.Range("D13").Value = 0.7741935483870968
$ws.Range("B14").Value = 0.6989247311827957
$ws.Range("C14").Value = 0.6989247311827957
$ws.Range("D14").Value = 0.6989247311827957
$ws.Range("E14").Value = 0.6989247311827957
$ws.Range("B15").Value = 0.7019927536231885
$ws.Range("C15").Value = 0.6600478468899521
$ws.Range("D15").Value = 0.6612903225806452
$ws.Range("B16").Value = 0.7008337229234846
$ws.Range("C16").Value = 0.6989247311827957
$ws.Range("D16").Value = 0.6819285466527922
$ws.Range("B17").Value = 0.6666666666666666
$ws.Range("C17").Value = 0.7368421052631579
$ws.Range("D17").Value = 0.7
$ws.Range("B18").Value = 0.803921568627451
$ws.Range("C18").Value = 0.7454545454545455
$ws.Range("D18").Value = 0.7735849056603775
$ws.Range("B19").Value = 0.7419354838709677
$ws.Range("C19").Value = 0.7419354838709677
$ws.Range("D19").Value = 0.7419354838709677
$ws.Range("E19").Value = 0.7419354838709677
$ws.Range("B20").Value = 0.7352941176470589
$ws.Range("C20").Value = 0.7411483253588517
$ws.Range("D20").Value = 0.7367924528301888
$ws.Range("B21").Value = 0.7478389205144423
$ws.Range("C21").Value = 0.7419354838709677
$ws.Range("D21").Value = 0.7435179549604384
$ws.Range("B22").Value = 0.6041666666666666
$ws.Range("C22").Value = 0.7631578947368421
$ws.Range("D22").Value = 0.6744186046511628
$ws.Range("B23").Value = 0.8
$ws.Range("C23").Value = 0.6545454545454545
$ws.Range("D23").Value = 0.7200000000000001
$ws.Range("B24").Value = 0.6989247311827957
$ws.Range("C24").Value = 0.6989247311827957
$ws.Range("D24").Value = 0.6989247311827957
$ws.Range("E24").Value = 0.6989247311827957
$ws.Range("B25").Value = 0.7020833333333334
$ws.Range("C25").Value = 0.7088516746411484
$ws.Range("D25").Value = 0.6972093023255814
$ws.Range("B26").Value = 0.7199820788530465
$ws.Range("C26").Value = 0.6989247311827957
$ws.Range("D26").Value = 0.7013753438359589
